$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '26.129.03'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  -2.34%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.572.38'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  -1.88%  '
$ws.Range('E4').Value = '  -0.55%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '208.76'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -1.53%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.498'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -3.14%  '
$ws.Range('E7').Value = '  -0.47%  '
$ws.Range('E8').Value = '  -1.79%  '
$ws.Range('E9').Value = '  -1.33%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '19.60'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -0.61%  '
$ws.Range('E11').Value = '  -0.19%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '1.791.86'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -1.93%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '1.592.12'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -0.66%  '
$ws.Range('B14').Value = 'Polkadot'
$ws.Range('C14').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '4.05'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -0.60%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.515'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -1.99%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '64.20'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  -1.38%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '26.108.27'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -2.29%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.0₃0723'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -2.28%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '7.32'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +2.45%  '
$ws.Range('B20').Value = 'BitcoinCash'
$ws.Range('C20').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '207.70'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -1.05%  '
$ws.Range('B21').Value = 'Dai'
$ws.Range('C21').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '1.00'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -0.42%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '4.24'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -1.60%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '2.19'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -2.03%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '8.82'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -2.98%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '143.72'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -0.06%  '
$ws.Range('E26').Value = '  -0.52%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '6.97'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -1.94%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '0.112'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -2.00%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '15.20'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -1.37%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '0.0505'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -0.61%  '
$ws.Range('E31').Value = '  -1.45%  '
$ws.Range('E32').Value = '  -2.22%  '
$ws.Range('E33').Value = '  +0.07%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '1.277.61'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -1.31%  '
$ws.Range('B35').Value = 'ImmutableX'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.611'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +2.08%  '
$ws.Range('B36').Value = 'HuobiToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '2.44'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -1.55%  '
$ws.Range('E37').Value = '  -1.61%  '
$ws.Range('E38').Value = '  -2.86%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '1.07'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -10.36%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.811'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -2.86%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '5.54'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +1.98%  '
$ws.Range('E42').Value = '  -3.05%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.762'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -2.73%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '62.28'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -1.35%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '1.704.83'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -2.02%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '89.03'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -1.60%  '
$ws.Range('B47').Value = 'RenderToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '1.51'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -3.28%  '
$ws.Range('B48').Value = 'BabyDogeCoin'
$ws.Range('C48').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.0₆0101'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -2.15%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.100'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -1.99%  '
$ws.Range('E50').Value = '  -1.73%  '
